$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap country names / data rows 215 and 216 (Islas Malvinas <-> Montserrat)
$ws.Range("A215").Value = "Montserrat"
$ws.Range("A216").Value = "Islas Malvinas"

$ws.Range("D215").Value = 12
$ws.Range("H215").Value = 1

$ws.Range("D216").Value = 13
$ws.Range("H216").Value = 0

# Row 4 - Estados Unidos
$ws.Range("D4").Value = 4849269
$ws.Range("E4").Value = 2573032

# Row 28 - Ucrania
$ws.Range("B28").Value = 230236
$ws.Range("C28").Value = 3774
$ws.Range("D28").Value = 101252
$ws.Range("E28").Value = 124554
$ws.Range("G28").Value = 33
$ws.Range("H28").Value = 4430

# Row 51 - Bielorrusia
$ws.Range("B51").Value = 80295
$ws.Range("E51").Value = 4290
$ws.Range("H51").Value = 857

# Row 63 - Armenia
$ws.Range("B63").Value = 52677
$ws.Range("C63").Value = 181
$ws.Range("D63").Value = 44710
$ws.Range("E63").Value = 6983
$ws.Range("G63").Value = 7
$ws.Range("H63").Value = 984

# Row 73 - Afganistan
$ws.Range("B73").Value = 39422
$ws.Range("C73").Value = 81
$ws.Range("D73").Value = 32879
$ws.Range("E73").Value = 5077
$ws.Range("G73").Value = 4
$ws.Range("H73").Value = 1466

# Row 79 - El Salvador
$ws.Range("E79").Value = 4289
$ws.Range("G79").Value = 2
$ws.Range("H79").Value = 865

# Row 81 - Australia
$ws.Range("B81").Value = 27149
$ws.Range("C81").Value = 13
$ws.Range("D81").Value = 24892
$ws.Range("E81").Value = 1363

# Row 113 - Georgia
$ws.Range("B113").Value = 8696
$ws.Range("C113").Value = 578
$ws.Range("D113").Value = 4619
$ws.Range("E113").Value = 4025
$ws.Range("G113").Value = 2
$ws.Range("H113").Value = 52

# Row 143 - Sri Lanka
$ws.Range("B143").Value = 3471
$ws.Range("C143").Value = 69
$ws.Range("E143").Value = 200

# Row 175 - Taiwan
$ws.Range("B175").Value = 518
$ws.Range("C175").Value = 1
$ws.Range("D175").Value = 485
